# Update examples and documentation:
#  - Reorder some of the variable lists in "Significant Components" and
#    "Included and Excluded" to a new canonical ordering.
#  - Refresh the computed factor-analysis results in "Loading Factors",
#    "All Refactor Variances" and "Final Variances" (row order + values)
#    to match the newly-ordered variable list.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Significant Components")
$ws.Range("C2").Value = '[''PPUNIT'' ''QNOHLTH'' ''QSERV'' ''QEXTRCT'' ''QESL'' ''QHISPC'' ''QEDLESHI'' ''QFHH''
 ''PERCAP'']'
$ws.Range("C4").Value = '[''MEDAGE'' ''QAGEDEP'' ''QSSBEN'']'
$ws.Range("C5").Value = '[''QRENTER'' ''QNOAUTO'' ''QPOVTY'']'
$ws.Range("C6").Value = '[''QAGEDEP'' ''QFEMALE'' ''QFEMLBR'']'

$ws = $wb.Worksheets.Item("Loading Factors")
$ws.Range("B2").Value = 0.7835523171479716
$ws.Range("C2").Value = -0.003789133447526956
$ws.Range("D2").Value = -0.1313207950294227
$ws.Range("E2").Value = -0.3772473370256653
$ws.Range("F2").Value = 0.1043526307693475
$ws.Range("A3").Value = 'QNOHLTH'
$ws.Range("B3").Value = 0.6484773982365736
$ws.Range("C3").Value = 0.4451428930904934
$ws.Range("D3").Value = -0.09382909758839381
$ws.Range("E3").Value = 0.3175828039673165
$ws.Range("F3").Value = -0.1542397762835904
$ws.Range("A4").Value = 'QSERV'
$ws.Range("B4").Value = 0.5277738002356493
$ws.Range("C4").Value = 0.3860111196285583
$ws.Range("D4").Value = -0.186875809908064
$ws.Range("E4").Value = 0.3543215291102447
$ws.Range("F4").Value = -0.03855995781613342
$ws.Range("A5").Value = 'QEXTRCT'
$ws.Range("B5").Value = 0.7523372382134341
$ws.Range("C5").Value = 0.1476097714664718
$ws.Range("D5").Value = -0.01125724324354274
$ws.Range("E5").Value = 0.09236455893403135
$ws.Range("F5").Value = -0.268847905743174
$ws.Range("A6").Value = 'QESL'
$ws.Range("B6").Value = 0.7907959801087903
$ws.Range("C6").Value = 0.1697374938532863
$ws.Range("D6").Value = -0.01961142838724014
$ws.Range("E6").Value = 0.2266196897906891
$ws.Range("F6").Value = -0.2859865792191026
$ws.Range("A7").Value = 'QHISPC'
$ws.Range("B7").Value = 0.8335417240452991
$ws.Range("C7").Value = 0.3477459695309182
$ws.Range("D7").Value = -0.1159860998568331
$ws.Range("E7").Value = 0.1412652495744246
$ws.Range("F7").Value = -0.09805697550088645
$ws.Range("A8").Value = 'QEDLESHI'
$ws.Range("B8").Value = 0.8755563542687488
$ws.Range("C8").Value = 0.2475053920687619
$ws.Range("D8").Value = 0.01400336361569062
$ws.Range("E8").Value = 0.2099414439935967
$ws.Range("F8").Value = -0.1400362615273635
$ws.Range("B9").Value = 0.5447678221922488
$ws.Range("C9").Value = 0.3097178103341896
$ws.Range("D9").Value = -0.08643488338332796
$ws.Range("E9").Value = 0.07397525862301159
$ws.Range("F9").Value = 0.2884044970735956
$ws.Range("B10").Value = 0.3668007650161663
$ws.Range("C10").Value = 0.8225015955555288
$ws.Range("D10").Value = -0.02259948686178994
$ws.Range("E10").Value = -0.03119507389113898
$ws.Range("F10").Value = -0.02428594150303115
$ws.Range("B11").Value = 0.4803068722416163
$ws.Range("C11").Value = 0.695654014931484
$ws.Range("D11").Value = -0.2531491196281418
$ws.Range("E11").Value = 0.2524414187005414
$ws.Range("F11").Value = 0.06893419244810839
$ws.Range("B12").Value = 0.1818961254000639
$ws.Range("C12").Value = 0.8639279902154168
$ws.Range("D12").Value = -0.1646723792240385
$ws.Range("E12").Value = 0.3092163416156365
$ws.Range("F12").Value = -0.01566262984896203
$ws.Range("A13").Value = 'QRENTER'
$ws.Range("B13").Value = -0.01907546856000766
$ws.Range("C13").Value = 0.2341572220789369
$ws.Range("D13").Value = -0.426819444523264
$ws.Range("E13").Value = 0.7484008130420402
$ws.Range("F13").Value = -0.1155640084771295
$ws.Range("A14").Value = 'QNOAUTO'
$ws.Range("B14").Value = 0.1274383130594738
$ws.Range("C14").Value = 0.08567042866239177
$ws.Range("D14").Value = -0.06014622622027892
$ws.Range("E14").Value = 0.6906315098292366
$ws.Range("F14").Value = 0.04382609482838686
$ws.Range("A15").Value = 'QPOVTY'
$ws.Range("B15").Value = 0.3961340619172994
$ws.Range("C15").Value = 0.1517075703026682
$ws.Range("D15").Value = -0.3133146679383343
$ws.Range("E15").Value = 0.5535940134509935
$ws.Range("F15").Value = 0.0922954857197491
$ws.Range("A16").Value = 'MEDAGE'
$ws.Range("B16").Value = -0.3172320250559871
$ws.Range("C16").Value = -0.242527512728358
$ws.Range("D16").Value = 0.7789709702313978
$ws.Range("E16").Value = -0.2970161777501997
$ws.Range("F16").Value = -0.05364863875777251
$ws.Range("A17").Value = 'QAGEDEP'
$ws.Range("B17").Value = -0.03648373506073081
$ws.Range("C17").Value = -0.1156683506291939
$ws.Range("D17").Value = 0.6860459936165838
$ws.Range("E17").Value = -0.09660788809046474
$ws.Range("F17").Value = 0.5943324668861167
$ws.Range("B18").Value = 0.01908852537982779
$ws.Range("C18").Value = -0.03712547503950558
$ws.Range("D18").Value = 0.7831494215062663
$ws.Range("E18").Value = -0.1342903125725542
$ws.Range("F18").Value = 0.0995559592412104
$ws.Range("A19").Value = 'QFEMALE'
$ws.Range("B19").Value = -0.05366633227337104
$ws.Range("C19").Value = -0.06387086926121399
$ws.Range("D19").Value = 0.184463602638826
$ws.Range("E19").Value = -0.006743986934936302
$ws.Range("F19").Value = 0.8303066316541788
$ws.Range("A20").Value = 'QFEMLBR'
$ws.Range("B20").Value = -0.2367556383835338
$ws.Range("C20").Value = 0.08379889872234192
$ws.Range("D20").Value = -0.03449523983892155
$ws.Range("E20").Value = 0.04713344724373027
$ws.Range("F20").Value = 0.8303066219177663

$ws = $wb.Worksheets.Item("All Refactor Variances")
$ws.Range("B2").Value = 5.207208521200397
$ws.Range("C2").Value = 2.879660807879201
$ws.Range("D2").Value = 2.320870209207491
$ws.Range("E2").Value = 2.280499951741452
$ws.Range("F2").Value = 1.987191277117984
$ws.Range("G2").Value = 1.6303292642431
$ws.Range("H2").Value = 1.07968958087668
$ws.Range("I2").Value = 4.722615406815523
$ws.Range("J2").Value = 3.480168177672616
$ws.Range("K2").Value = 2.298524093778104
$ws.Range("L2").Value = 2.143689086244914
$ws.Range("M2").Value = 2.01807390787403
$ws.Range("N2").Value = 4.995669856720731
$ws.Range("O2").Value = 2.752559874573491
$ws.Range("P2").Value = 2.184268141804553
$ws.Range("Q2").Value = 2.119818616378288
$ws.Range("R2").Value = 2.038891349445035
$ws.Range("B3").Value = 0.192859574859274
$ws.Range("C3").Value = 0.106654103995526
$ws.Range("D3").Value = 0.08595815589657373
$ws.Range("E3").Value = 0.08446296117560932
$ws.Range("F3").Value = 0.07359967693029569
$ws.Range("G3").Value = 0.06038256534233703
$ws.Range("H3").Value = 0.0399885029954326
$ws.Range("I3").Value = 0.2248864479435963
$ws.Range("J3").Value = 0.1657222941748865
$ws.Range("K3").Value = 0.1094535282751478
$ws.Range("L3").Value = 0.1020804326783293
$ws.Range("M3").Value = 0.09609875751781095
$ws.Range("N3").Value = 0.2629299924589859
$ws.Range("O3").Value = 0.1448715723459732
$ws.Range("P3").Value = 0.1149614811476081
$ws.Range("Q3").Value = 0.1115694008620152
$ws.Range("R3").Value = 0.1073100710234229
$ws.Range("B4").Value = 0.192859574859274
$ws.Range("C4").Value = 0.2995136788548
$ws.Range("D4").Value = 0.3854718347513737
$ws.Range("E4").Value = 0.469934795926983
$ws.Range("F4").Value = 0.5435344728572787
$ws.Range("G4").Value = 0.6039170381996157
$ws.Range("H4").Value = 0.6439055411950483
$ws.Range("I4").Value = 0.2248864479435963
$ws.Range("J4").Value = 0.3906087421184828
$ws.Range("K4").Value = 0.5000622703936306
$ws.Range("L4").Value = 0.6021427030719598
$ws.Range("M4").Value = 0.6982414605897708
$ws.Range("N4").Value = 0.2629299924589859
$ws.Range("O4").Value = 0.407801564804959
$ws.Range("P4").Value = 0.5227630459525672
$ws.Range("Q4").Value = 0.6343324468145823
$ws.Range("R4").Value = 0.7416425178380053
$ws.Range("B5").Value = 0.2995153209915521
$ws.Range("C5").Value = 0.1656362574510271
$ws.Range("D5").Value = 0.1334949777525455
$ws.Range("E5").Value = 0.1311729062291518
$ws.Range("F5").Value = 0.1143019778859168
$ws.Range("G5").Value = 0.09377550196302205
$ws.Range("H5").Value = 0.06210305772678466
$ws.Range("I5").Value = 0.322075471934373
$ws.Range("J5").Value = 0.2373423858773853
$ws.Range("K5").Value = 0.156755985504925
$ws.Range("L5").Value = 0.1461964641745949
$ws.Range("M5").Value = 0.1376296925087218
$ws.Range("N5").Value = 0.3545238927582856
$ws.Range("O5").Value = 0.1953388173702537
$ws.Range("P5").Value = 0.1550092913803504
$ws.Range("Q5").Value = 0.1504355510620616
$ws.Range("R5").Value = 0.1446924474290487

$ws = $wb.Worksheets.Item("Final Variances")
$ws.Range("B2").Value = 4.995669856720731
$ws.Range("C2").Value = 2.752559874573491
$ws.Range("D2").Value = 2.184268141804553
$ws.Range("E2").Value = 2.119818616378288
$ws.Range("F2").Value = 2.038891349445035
$ws.Range("B3").Value = 0.2629299924589859
$ws.Range("C3").Value = 0.1448715723459732
$ws.Range("D3").Value = 0.1149614811476081
$ws.Range("E3").Value = 0.1115694008620152
$ws.Range("F3").Value = 0.1073100710234229
$ws.Range("B4").Value = 0.2629299924589859
$ws.Range("C4").Value = 0.407801564804959
$ws.Range("D4").Value = 0.5227630459525672
$ws.Range("E4").Value = 0.6343324468145823
$ws.Range("F4").Value = 0.7416425178380053
$ws.Range("B5").Value = 0.3545238927582856
$ws.Range("C5").Value = 0.1953388173702537
$ws.Range("D5").Value = 0.1550092913803504
$ws.Range("E5").Value = 0.1504355510620616
$ws.Range("F5").Value = 0.1446924474290487

$ws = $wb.Worksheets.Item("Included and Excluded")
$ws.Range("B2").Value = '[[''PPUNIT'', ''QNOHLTH'', ''QSERV'', ''QEXTRCT'', ''QESL'', ''QHISPC'', ''QEDLESHI'', ''QFHH'', ''PERCAP'', ''MDHSEVAL'', ''QRICH'', ''MEDAGE'', ''QAGEDEP'', ''QSSBEN'', ''QRENTER'', ''QNOAUTO'', ''QPOVTY'', ''QFEMALE'', ''QFEMLBR'']]'
